$d = $word.ActiveDocument

# Locate the Specific Aims intro paragraph by its distinctive opening text,
# rather than relying on a hard-coded paragraph index.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Alternative splicing of mutually exclusive exons")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the target paragraph to rewrite."
}

$r = $target.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "Alternative splicing is an important factor in gene regulation; however, its regulation is still incompletely understood. While splicing decisions are driven by information in the local sequence context, outcomes can be correlated over large distances. A particularly common example of this is mutually exclusive exons. These occur when exons in a cluster are anticorrelated so only one is included in the mature transcript. Several mechanisms for producing this behavior have been identified, including spliceosome incompatibility, 5’-splice site/branch point steric occlusion, and RNA secondary structures but most mutually exclusive exon clusters in humans still have no known mechanism. While these long-distance anti-correlations present challenges as a study system because they require measuring multiple exons simultaneously in a single transcript. Work on the function of mutually exclusive exons has borne fruit. A common theme is switch-like changes in protein function that occur in cell differentiation or development. Pyruvate Kinase M (PKM), ketohexokinase, and CaV1.2 are some of the better-known examples of important regulatory events that occur through mutually exclusive exon switches. PKM is converted to a constitutively active form during differentiation, reversion of this change in cancer cells is a major driver of the Warburg effect. Ketohexokinase can be switched between high and low affinity forms to control fructose metabolism, aberrant expression of the high affinity form causes pathological cardiac hypertrophy while the low affinity form is a major driver of hepatocellular carcinoma. Finally, mutations that shift the isoform distribution of CaV1.2 cause timothy syndrome characterized by catastrophic developmental, neurological, and cardiac symptoms with a life expectancy of ~2.5 years. Understanding mutually exclusive splicing will not only deepen our knowledge of development and differentiation but bring us closer to treatments for a wide variety of human diseases. "
